$wb = $excel.ActiveWorkbook
$wsUserStories = $wb.Worksheets.Item("User Stories")
$wsIteration1  = $wb.Worksheets.Item("StoriesDetailsIteration1")

# --- Cell value edits on "StoriesDetailsIteration1" (sheet2) ---
# Order below matters: new shared-string entries are appended in first-use
# order, and the target file expects them in this sequence:
#   146 Wyatt, 147 Chris, 148 David, 149 Waleed, 150 Wyatt/David

# J9: Wyatt (new cell, no style)
$wsIteration1.Range("J9").Value = "Wyatt"

# J4: Chris (new cell, needs style 17 - vertical top alignment)
$wsIteration1.Range("J4").VerticalAlignment = -4160
$wsIteration1.Range("J4").Value = "Chris"

# J7: David (new cell, no style)
$wsIteration1.Range("J7").Value = "David"

# Row 5: remove the old "Done" note from K5 and add "Waleed" note to J5 instead
$wsIteration1.Range("K5").ClearContents()
$wsIteration1.Range("J5").Value = "Waleed"

# J20: Wyatt/David (existing empty cell already styled s="17")
$wsIteration1.Range("J20").Value = "Wyatt/David"

# J6: Waleed (new cell, no style) - reuses shared string created at J5
$wsIteration1.Range("J6").Value = "Waleed"

# K9: a single space note (reuses existing shared string)
$wsIteration1.Range("K9").Value = " "

# J14, J21, J27: Chris (existing empty cells already styled s="17")
$wsIteration1.Range("J14").Value = "Chris"
$wsIteration1.Range("J21").Value = "Chris"
$wsIteration1.Range("J27").Value = "Chris"

# K28: Done (existing empty cell already styled s="17")
$wsIteration1.Range("K28").Value = "Done"

# --- Sheet view / selection changes ---
# "User Stories" (sheet1): selection moves to B3, scroll resets (no topLeftCell)
$null = $wsUserStories.Range("B3").Select()

# Re-activate "StoriesDetailsIteration1" so it stays the tab-selected sheet
# (it was the active sheet before these edits), then update its selection.
$null = $wsIteration1.Activate()
$null = $wsIteration1.Range("B27:I27").Select()
